$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1624.5
$ws.Range("J17").Value = 1624.5
$ws.Range("L17").Value = 4873.5
$ws.Range("N17").Value = -5209.5

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2273.75
$ws.Range("I40").Value = 1598
$ws.Range("J40").Value = 2499
$ws.Range("K40").Value = 1598
$ws.Range("L40").Value = 2499
$ws.Range("M40").Value = -1423
$ws.Range("N40").Value = -2849

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 500
$ws.Range("I125").Value = 500
$ws.Range("K125").Value = 4500
$ws.Range("M125").Value = -2040

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3163.8823
$ws.Range("I137").Value = 2513.2856
$ws.Range("J137").Value = 6200
$ws.Range("K137").Value = 7539.8568
$ws.Range("L137").Value = 18600
$ws.Range("M137").Value = -4989.8568
$ws.Range("N137").Value = -23700

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 144
$ws.Range("I5").Value = 136.8
$ws.Range("K5").Value = 136.8
$ws.Range("M5").Value = -24.80000000000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4790.1816
$ws.Range("I102").Value = 4790.1816
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4790.1816
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -3168.1816

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3750
$ws.Range("I132").Value = 3750
$ws.Range("K132").Value = 11250
$ws.Range("M132").Value = -8720

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 144
$ws.Range("I4").Value = 136.8
$ws.Range("K4").Value = 136.8
$ws.Range("M4").Value = -21.80000000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 708.5833
$ws.Range("I99").Value = 708.5833
$ws.Range("K99").Value = 708.5833
$ws.Range("M99").Value = 789.4167

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 599.5
$ws.Range("I107").Value = 599.5
$ws.Range("K107").Value = 599.5
$ws.Range("M107").Value = 1320.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.333336
$ws.Range("I7").Value = 109.333336
$ws.Range("K7").Value = 109.333336
$ws.Range("M7").Value = 3.666663999999997

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 399.72726
$ws.Range("I22").Value = 418
$ws.Range("J22").Value = 351
$ws.Range("K22").Value = 418
$ws.Range("L22").Value = 351
$ws.Range("M22").Value = -68
$ws.Range("N22").Value = -1051

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2684.7334
$ws.Range("I31").Value = 1674.25
$ws.Range("K31").Value = 1674.25
$ws.Range("M31").Value = -1379.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1572.2858
$ws.Range("I33").Value = 1584.3334
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 1584.3334
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -1205.3334
$ws.Range("N33").Value = -2258

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2684.7334
$ws.Range("I34").Value = 1674.25
$ws.Range("K34").Value = 1674.25
$ws.Range("M34").Value = -1472.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4800
$ws.Range("I36").Value = 4800
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4800
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -4412

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 4800
$ws.Range("I40").Value = 4800
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4800
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4640

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 955.8823
$ws.Range("I4").Value = 395.8
$ws.Range("J4").Value = 1756
$ws.Range("K4").Value = 1187.4
$ws.Range("L4").Value = 5268
$ws.Range("M4").Value = -1075.4
$ws.Range("N4").Value = -5492

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1592.9445
$ws.Range("I121").Value = 1050
$ws.Range("J121").Value = 1660.8125
$ws.Range("K121").Value = 3150
$ws.Range("L121").Value = 4982.4375
$ws.Range("M121").Value = -1840
$ws.Range("N121").Value = -7602.4375

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 55.333332
$ws.Range("I2").Value = 22.666666
$ws.Range("J2").Value = 88
$ws.Range("K2").Value = 22.666666
$ws.Range("L2").Value = 88
$ws.Range("M2").Value = 90.33333400000001
$ws.Range("N2").Value = -314

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4000000
$ws.Range("J7").Value = 4000000
$ws.Range("L7").Value = 4000000
$ws.Range("N7").Value = -4000224

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 4000000
$ws.Range("J8").Value = 4000000
$ws.Range("L8").Value = 4000000
$ws.Range("N8").Value = -4000278

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 9500.5
$ws.Range("I24").Value = 9000
$ws.Range("J24").Value = 10001
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10001
$ws.Range("M24").Value = -8827
$ws.Range("N24").Value = -10347

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15046
$ws.Range("I46").Value = 15046
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 15046
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -14890

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5941.1763
$ws.Range("I16").Value = 4281.25
$ws.Range("K16").Value = 4281.25
$ws.Range("M16").Value = -4111.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 534.625
$ws.Range("I22").Value = 454.66666
$ws.Range("K22").Value = 454.66666
$ws.Range("M22").Value = -159.66666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 534.625
$ws.Range("I27").Value = 454.66666
$ws.Range("K27").Value = 454.66666
$ws.Range("M27").Value = -347.66666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 649.5
$ws.Range("I46").Value = 599.3333
$ws.Range("K46").Value = 599.3333
$ws.Range("M46").Value = -411.3333

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 47941.816
$ws.Range("I136").Value = 4741.1
$ws.Range("K136").Value = 14223.3
$ws.Range("M136").Value = -11673.3

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4251.4614
$ws.Range("I136").Value = 4251.4614
$ws.Range("K136").Value = 12754.3842
$ws.Range("M136").Value = -10204.3842
